# Update the cryptos price (D) and 1h-volume-change (E) columns.
# D-column values are plain text in the source workbook (e.g. "39.551.97",
# "226.44"). Assigning a bare numeric-looking string via .Value would make
# Excel auto-convert it to a Number, which would not match the original
# "inline string" cell type. Prefixing with an apostrophe forces Excel to
# keep it as text (quotePrefix); re-applying the "Normal" style afterwards
# clears the visual quote-prefix styling so the cell's style stays
# identical to before (no explicit style index), matching the source file.
#
# E-column values already contain non-numeric characters (leading/trailing
# spaces, a trailing "%") so Excel keeps them as plain text without any
# extra trick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
  $ws.Range($addr).Value = "'" + $text
  $ws.Range($addr).Style = "Normal"
}

$updates = @(
  @{ Row = 2;  D = "39.551.97";  E = "  +1.91%  " }
  @{ Row = 3;  D = "2.153.90";   E = "  +2.28%  " }
  @{ Row = 4;             E = "  -0.02%  " }
  @{ Row = 5;  D = "226.44";     E = "  -0.69%  " }
  @{ Row = 6;  D = "0.618";      E = "  +0.33%  " }
  @{ Row = 7;  D = "62.64";      E = "  +1.09%  " }
  @{ Row = 8;             E = "  -0.01%  " }
  @{ Row = 9;             E = "  -0.23%  " }
  @{ Row = 10; D = "0.0840";     E = "  -0.37%  " }
  @{ Row = 11;            E = "  -0.05%  " }
  @{ Row = 12; D = "15.79";      E = "  -1.16%  " }
  @{ Row = 13; D = "2.473.79";   E = "  +2.37%  " }
  @{ Row = 14; D = "21.65";      E = "  -1.69%  " }
  @{ Row = 15;            E = "  -0.01%  " }
  @{ Row = 16; D = "5.46" }
  @{ Row = 17; D = "2.166.66";   E = "  +3.89%  " }
  @{ Row = 18; D = "39.532.64";  E = "  +1.56%  " }
  @{ Row = 19;            E = "  -0.21%  " }
  @{ Row = 20;            E = "  -0.05%  " }
  @{ Row = 21; D = "0.0₃0850";   E = "  +0.62%  " }
  @{ Row = 22; D = "227.37";     E = "  +0.05%  " }
  @{ Row = 23;            E = "  +0.02%  " }
  @{ Row = 24; D = "2.35";       E = "  +1.35%  " }
  @{ Row = 25;            E = "  -0.80%  " }
  @{ Row = 26; D = "170.57";     E = "  -0.12%  " }
  @{ Row = 27;            E = "  -2.81%  " }
  @{ Row = 28;            E = "  +1.31%  " }
  @{ Row = 29;            E = "  +0.71%  " }
  @{ Row = 30;            E = "  +1.09%  " }
  @{ Row = 31; D = "2.69";       E = "  +4.91%  " }
  @{ Row = 32;            E = "  +0.37%  " }
  @{ Row = 33;            E = "  -0.29%  " }
  @{ Row = 34; D = "4.70";       E = "  -1.79%  " }
  @{ Row = 35; D = "6.96";       E = "  -2.99%  " }
  @{ Row = 36;            E = "  +0.03%  " }
  @{ Row = 37; D = "3.80";       E = "  +8.35%  " }
  @{ Row = 38;            E = "  +0.87%  " }
  @{ Row = 39;            E = "  -0.06%  " }
  @{ Row = 40; D = "4.87";       E = "  +17.70%  " }
  @{ Row = 41; D = "102.49";     E = "  +0.38%  " }
  @{ Row = 42;            E = "  -1.21%  " }
  @{ Row = 43; D = "17.63";      E = "  -2.34%  " }
  @{ Row = 44; D = "1.514.25";   E = "  -0.98%  " }
  @{ Row = 45;            E = "  +0.04%  " }
  @{ Row = 46; D = "7.82";       E = "  +0.88%  " }
  @{ Row = 47;            E = "  +0.13%  " }
  @{ Row = 48; D = "0.0919";     E = "  +0.17%  " }
  @{ Row = 49;            E = "  -0.14%  " }
  @{ Row = 50;            E = "  +0.79%  " }
  @{ Row = 51; D = "2.357.23";   E = "  +2.27%  " }
)

foreach ($u in $updates) {
  $row = $u.Row
  if ($u.ContainsKey("D")) {
    Set-TextCell "D$row" $u["D"]
  }
  if ($u.ContainsKey("E")) {
    $ws.Range("E$row").Value = $u["E"]
  }
}
